$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D..K shift to F..M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from the (now shifted) original column F onto the
# two freshly inserted columns D:E, per contiguous data block, so the new
# cells pick up the same style id as the rest of the row (dates vs numbers)
# without generating brand-new style entries.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(8,4).Value = 639100
$ws.Cells.Item(8,5).Value = 641100
$ws.Cells.Item(9,4).Value = 222000
$ws.Cells.Item(9,5).Value = 228400
$ws.Cells.Item(10,4).Value = 417100
$ws.Cells.Item(10,5).Value = 412700
$ws.Cells.Item(12,4).Value = "NA"
$ws.Cells.Item(12,5).Value = "NA"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = 2800
$ws.Cells.Item(14,5).Value = 700
$ws.Cells.Item(15,4).Value = 146900
$ws.Cells.Item(15,5).Value = 167800
$ws.Cells.Item(17,4).Value = 494400
$ws.Cells.Item(17,5).Value = 518300
$ws.Cells.Item(18,4).Value = 144700
$ws.Cells.Item(18,5).Value = 122800
$ws.Cells.Item(20,4).Value = -8000
$ws.Cells.Item(20,5).Value = 2000
$ws.Cells.Item(21,4).Value = 283600
$ws.Cells.Item(21,5).Value = 292600
$ws.Cells.Item(22,4).Value = 84000
$ws.Cells.Item(22,5).Value = 82200
$ws.Cells.Item(23,4).Value = 52700
$ws.Cells.Item(23,5).Value = 42600
$ws.Cells.Item(24,4).Value = 22900
$ws.Cells.Item(24,5).Value = 12900
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = 29800
$ws.Cells.Item(26,5).Value = 29700
$ws.Cells.Item(27,4).Value = 29800
$ws.Cells.Item(27,5).Value = 29700
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = 400
$ws.Cells.Item(29,5).Value = -7600
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = 8000
$ws.Cells.Item(32,5).Value = -2000
$ws.Cells.Item(33,4).Value = 30200
$ws.Cells.Item(33,5).Value = 22100
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = 30200
$ws.Cells.Item(35,5).Value = 22100
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(41,4).Value = 176400
$ws.Cells.Item(41,5).Value = 353900
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(43,4).Value = 202800
$ws.Cells.Item(43,5).Value = 227200
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,5).Value = 0
$ws.Cells.Item(45,4).Value = 100800
$ws.Cells.Item(45,5).Value = 99800
$ws.Cells.Item(46,4).Value = 480000
$ws.Cells.Item(46,5).Value = 680900
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(48,4).Value = 5582300
$ws.Cells.Item(48,5).Value = 5524700
$ws.Cells.Item(49,4).Value = 2870400
$ws.Cells.Item(49,5).Value = 2902700
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 204700
$ws.Cells.Item(52,5).Value = 206600
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 9137400
$ws.Cells.Item(54,5).Value = 9314900
$ws.Cells.Item(57,4).Value = 48200
$ws.Cells.Item(57,5).Value = 41900
$ws.Cells.Item(58,4).Value = 21000
$ws.Cells.Item(58,5).Value = 15500
$ws.Cells.Item(59,4).Value = 493300
$ws.Cells.Item(59,5).Value = 568500
$ws.Cells.Item(60,4).Value = 562500
$ws.Cells.Item(60,5).Value = 625900
$ws.Cells.Item(61,4).Value = 6068500
$ws.Cells.Item(61,5).Value = 5834200
$ws.Cells.Item(62,4).Value = 1324600
$ws.Cells.Item(62,5).Value = 1300000
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 7955600
$ws.Cells.Item(66,5).Value = 7760100
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = -313700
$ws.Cells.Item(72,5).Value = -343900
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 1181800
$ws.Cells.Item(76,5).Value = 1554800
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(81,4).Value = 30200
$ws.Cells.Item(81,5).Value = 22100
$ws.Cells.Item(83,4).Value = 146900
$ws.Cells.Item(83,5).Value = 167800
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = 230600
$ws.Cells.Item(89,5).Value = 241800
$ws.Cells.Item(91,4).Value = -202200
$ws.Cells.Item(91,5).Value = -182500
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = -202200
$ws.Cells.Item(94,5).Value = -143500
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = -196900
$ws.Cells.Item(100,5).Value = -6700
$ws.Cells.Item(101,4).Value = -9000
$ws.Cells.Item(101,5).Value = 2200
$ws.Cells.Item(102,4).Value = -177500
$ws.Cells.Item(102,5).Value = 93800

# Fix up the handful of cells in the shifted original column (now F)
# whose values were also updated in this data refresh.
$ws.Cells.Item(45,6).Value = 145600
$ws.Cells.Item(46,6).Value = 637900
$ws.Cells.Item(48,6).Value = 5427600
$ws.Cells.Item(52,6).Value = 213200
$ws.Cells.Item(54,6).Value = 9209900
$ws.Cells.Item(59,6).Value = 553900
$ws.Cells.Item(60,6).Value = 616700
$ws.Cells.Item(62,6).Value = 1281200
$ws.Cells.Item(66,6).Value = 7709600
$ws.Cells.Item(72,6).Value = -366000
$ws.Cells.Item(76,6).Value = 1500300
